$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "Custom Ref" column header and its value for row 2
$ws.Range("Q1").Value = "Custom Ref"
$ws.Range("Q2").Value = "Test5"

# Move the active selection to Q3 (as if Enter was pressed after typing into Q2)
[void]$ws.Range("Q3").Select()
